$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.253.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.88%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.580.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.78%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.24%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.75%  '

# Row 7
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.23%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.586.36'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.28%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.62%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.71%  '

# Row 12
$ws.Range("E12").Value = '  +11.53%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.345'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.89%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.034.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.84%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.221.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.80%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.72%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.20%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.584.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.98%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '337.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.65%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.52%  '

# Row 23
$ws.Range("E23").Value = '  +0.10%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.49%  '

# Row 25
$ws.Range("E25").Value = '  +5.19%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.01%  '

# Row 27
$ws.Range("E27").Value = '  -2.72%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.17%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0777'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.92%  '

# Row 30
$ws.Range("E30").Value = '  +0.03%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '161.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.24%  '

# Row 32
$ws.Range("E32").Value = '  -1.67%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.07'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.06%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.02'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.11%  '

# Row 36
$ws.Range("E36").Value = '  -0.44%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.882'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.49%  '

# Row 38
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.876'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.34%  '

# Row 39
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.48'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.05%  '

# Row 40
$ws.Range("E40").Value = '  -1.69%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '295.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.42%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.10%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '131.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.66%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0973'
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.596'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.15%  '

# Row 47
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0535'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.05%  '

# Row 48
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.16%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.17%  '

# Row 50
$ws.Range("E50").Value = '  -2.26%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.44'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.39%  '
